# Apply the "altered emx documentation, mref added to advanced example (data)" edit.
#
# Summary of the change:
#  - patients!A2:A4 display names get an underscore-joined form instead of a
#    space-joined one ("john doe" -> "john_doe", "jane doe" -> "jane_doe"),
#    and the "pape doe" typo is fixed to "papa_doe".
#  - The view/selection state flips: the "patients" sheet becomes the active
#    tab (selection E7) while "attributes" is no longer the active tab
#    (selection moves to C3).

$wb  = $excel.ActiveWorkbook
$wsPatients   = $wb.Worksheets.Item("patients")
$wsAttributes = $wb.Worksheets.Item("attributes")

# --- data edits: patients sheet, column A (displayName) ---
$wsPatients.Range("A2").Value = "john_doe"
$wsPatients.Range("A3").Value = "jane_doe"
$wsPatients.Range("A4").Value = "papa_doe"

# --- view/selection edits ---
# Select attributes!C3 first (while it's still the active sheet) so its
# stored selection updates, then make patients the active tab with E7
# selected, matching the target sheetViews state.
$wsAttributes.Range("C3").Select()
$wsPatients.Activate()
$wsPatients.Range("E7").Select()
